# chore: update Sheets via scheduled runner
#
# Refreshes cached market-board figures (currentAveragePrice / NQ / HQ)
# and the derived Leve profit columns (H:N) for the rows whose source
# prices moved since the last pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1171.25
$ws.Range("I18").Value = 307.4
$ws.Range("J18").Value = 2611
$ws.Range("K18").Value = 307.4
$ws.Range("L18").Value = 2611
$ws.Range("M18").Value = -23.39999999999998
$ws.Range("N18").Value = -3179

$ws.Range("H132").Value = 15526.586
$ws.Range("I132").Value = 16462.477
$ws.Range("J132").Value = 3360
$ws.Range("K132").Value = 49387.431
$ws.Range("L132").Value = 10080
$ws.Range("M132").Value = -46857.431
$ws.Range("N132").Value = -15140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25649.79
$ws.Range("I32").Value = 5776.8086
$ws.Range("J32").Value = 119052.8
$ws.Range("K32").Value = 5776.8086
$ws.Range("L32").Value = 119052.8
$ws.Range("M32").Value = -5489.8086
$ws.Range("N32").Value = -119626.8

$ws.Range("H41").Value = 17318.666
$ws.Range("I41").Value = 978
$ws.Range("K41").Value = 978
$ws.Range("M41").Value = -564

$ws.Range("H45").Value = 1192.8572
$ws.Range("I45").Value = 1191.6666
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1191.6666
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -814.6666
$ws.Range("N45").Value = -1954

$ws.Range("H61").Value = 2940.3635
$ws.Range("I61").Value = 2364.5625
$ws.Range("K61").Value = 2364.5625
$ws.Range("M61").Value = -2152.5625

$ws.Range("H74").Value = 3174.5173
$ws.Range("I74").Value = 967.87177
$ws.Range("J74").Value = 7703.9473
$ws.Range("K74").Value = 967.87177
$ws.Range("L74").Value = 7703.9473
$ws.Range("M74").Value = -93.87176999999997
$ws.Range("N74").Value = -9451.9473

$ws.Range("H77").Value = 3174.5173
$ws.Range("I77").Value = 967.87177
$ws.Range("J77").Value = 7703.9473
$ws.Range("K77").Value = 4839.35885
$ws.Range("L77").Value = 38519.7365
$ws.Range("M77").Value = -471.3588499999996
$ws.Range("N77").Value = -47255.7365

$ws.Range("H132").Value = 2211.9827
$ws.Range("I132").Value = 1715.1428
$ws.Range("K132").Value = 5145.428400000001
$ws.Range("M132").Value = -2615.428400000001

$ws.Range("H136").Value = 2940.3635
$ws.Range("I136").Value = 2364.5625
$ws.Range("K136").Value = 7093.6875
$ws.Range("M136").Value = -4543.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1297.7307
$ws.Range("I20").Value = 1074.0625
$ws.Range("J20").Value = 1655.6
$ws.Range("K20").Value = 1074.0625
$ws.Range("L20").Value = 1655.6
$ws.Range("M20").Value = -827.0625
$ws.Range("N20").Value = -2149.6

$ws.Range("H61").Value = 6797
$ws.Range("I61").Value = 6797
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6797
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6484
$ws.Range("N61").ClearContents()

$ws.Range("H86").Value = 1281.0952
$ws.Range("I86").Value = 1193.0714
$ws.Range("J86").Value = 1457.1428
$ws.Range("K86").Value = 1193.0714
$ws.Range("L86").Value = 1457.1428
$ws.Range("M86").Value = -70.07140000000004
$ws.Range("N86").Value = -3703.1428

$ws.Range("H89").Value = 1281.0952
$ws.Range("I89").Value = 1193.0714
$ws.Range("J89").Value = 1457.1428
$ws.Range("K89").Value = 5965.357
$ws.Range("L89").Value = 7285.714
$ws.Range("M89").Value = -349.357
$ws.Range("N89").Value = -18517.714

$ws.Range("H94").Value = 747.4815
$ws.Range("I94").Value = 601.64703
$ws.Range("J94").Value = 995.4
$ws.Range("K94").Value = 601.64703
$ws.Range("L94").Value = 995.4
$ws.Range("M94").Value = -150.64703
$ws.Range("N94").Value = -1897.4

$ws.Range("H99").Value = 2087.375
$ws.Range("I99").Value = 1783.1666
$ws.Range("K99").Value = 1783.1666
$ws.Range("M99").Value = -285.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3306.0217
$ws.Range("I31").Value = 1401.6451
$ws.Range("J31").Value = 7241.7334
$ws.Range("K31").Value = 1401.6451
$ws.Range("L31").Value = 7241.7334
$ws.Range("M31").Value = -1106.6451
$ws.Range("N31").Value = -7831.7334

$ws.Range("H34").Value = 3306.0217
$ws.Range("I34").Value = 1401.6451
$ws.Range("J34").Value = 7241.7334
$ws.Range("K34").Value = 1401.6451
$ws.Range("L34").Value = 7241.7334
$ws.Range("M34").Value = -1199.6451
$ws.Range("N34").Value = -7645.7334

$ws.Range("H122").Value = 1708
$ws.Range("I122").Value = 1473.3334
$ws.Range("J122").Value = 2764
$ws.Range("K122").Value = 4420.0002
$ws.Range("L122").Value = 8292
$ws.Range("M122").Value = -1970.0002
$ws.Range("N122").Value = -13192

$ws.Range("H132").Value = 1974.2084
$ws.Range("I132").Value = 1420.3
$ws.Range("J132").Value = 4743.75
$ws.Range("K132").Value = 4260.9
$ws.Range("L132").Value = 14231.25
$ws.Range("M132").Value = -1730.9
$ws.Range("N132").Value = -19291.25

$ws.Range("H141").Value = 353585.4
$ws.Range("I141").Value = 41000
$ws.Range("J141").Value = 382002.28
$ws.Range("K141").Value = 41000
$ws.Range("L141").Value = 382002.28
$ws.Range("M141").Value = -35820
$ws.Range("N141").Value = -392362.28

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3365.7058
$ws.Range("I7").Value = 2937.3333
$ws.Range("J7").Value = 3457.5
$ws.Range("K7").Value = 2937.3333
$ws.Range("L7").Value = 3457.5
$ws.Range("M7").Value = -2825.3333
$ws.Range("N7").Value = -3681.5

$ws.Range("H40").Value = 3300
$ws.Range("I40").Value = 2300
$ws.Range("J40").Value = 3476.4707
$ws.Range("K40").Value = 2300
$ws.Range("L40").Value = 3476.4707
$ws.Range("M40").Value = -2164
$ws.Range("N40").Value = -3748.4707

$ws.Range("H122").Value = 3979.4736
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4033.889
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 12101.667
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -17001.667

$ws.Range("H126").Value = 3365.7058
$ws.Range("I126").Value = 2937.3333
$ws.Range("J126").Value = 3457.5
$ws.Range("K126").Value = 8811.999899999999
$ws.Range("L126").Value = 10372.5
$ws.Range("M126").Value = -6341.999899999999
$ws.Range("N126").Value = -15312.5

$ws.Range("H132").Value = 4100.778
$ws.Range("I132").Value = 2719.125
$ws.Range("K132").Value = 8157.375
$ws.Range("M132").Value = -5627.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 880
$ws.Range("I100").Value = 966.6667
$ws.Range("J100").Value = 750
$ws.Range("K100").Value = 1933.3334
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -1392.3334
$ws.Range("N100").Value = -2582

$ws.Range("H122").Value = 1583.25
$ws.Range("I122").Value = 1327.421
$ws.Range("J122").Value = 2123.3333
$ws.Range("K122").Value = 3982.263
$ws.Range("L122").Value = 6369.999899999999
$ws.Range("M122").Value = -1532.263
$ws.Range("N122").Value = -11269.9999

$ws.Range("H126").Value = 91939.27
$ws.Range("I126").Value = 125485.375
$ws.Range("J126").Value = 2483
$ws.Range("K126").Value = 376456.125
$ws.Range("L126").Value = 7449
$ws.Range("M126").Value = -373986.125
$ws.Range("N126").Value = -12389

$ws.Range("H136").Value = 2958.66
$ws.Range("I136").Value = 1284.9678
$ws.Range("J136").Value = 5689.421
$ws.Range("K136").Value = 3854.9034
$ws.Range("L136").Value = 17068.263
$ws.Range("M136").Value = -1304.9034
$ws.Range("N136").Value = -22168.263
